# Applies the "Actualizacion del repositorio global ... 2024-12-12" edit:
#  - A1 gets the label "Encuestado" (previously blank, formatted with style s=1)
#  - Row 14 ("Nota:" / Likert scale note, merged B14:Q14) is removed entirely
#  - Two new summary rows are appended: MIN and MAX of the tutor-age column (B2:B13)
#  - A handful of data rows lose their manual/custom row height (revert to auto)
#  - The active selection moves to B17 (just under the new MAX formula)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label the first column header (was an empty, styled cell).
$ws.Range("A1").Value = "Encuestado"

# Drop the old "Nota:" / Likert-scale explanation row (14) completely,
# including its B14:Q14 merge.
$ws.Rows(14).EntireRow.Delete()

# Add the MIN/MAX summary formulas directly below the data table.
$ws.Range("B15").Formula = "=MIN(B2:B13)"
$ws.Range("B16").Formula = "=MAX(B2:B13)"

# These rows return to automatic row height (their manual height is cleared).
$ws.Rows("4").AutoFit()
$ws.Rows("9").AutoFit()
$ws.Rows("10").AutoFit()
$ws.Rows("11").AutoFit()

# Leave the selection on B17, right below the new formulas.
$ws.Activate() | Out-Null
$ws.Range("B17").Select() | Out-Null
